# Applies: "Add styles to the new paragraphs" (see #4)
#
# 1. Defines three new character styles (GaNStyle, GaNParagraph, GaNLinks)
#    in styles.xml, appended after the existing "UnresolvedMention" style.
# 2. Applies GaNStyle to the run(s) containing the "Informace v této
#    příručce ..." paragraph (4 occurrences), while also trimming the
#    duplicated trailing date-range text that had accidentally been left
#    in the run.
# 3. Applies GaNLinks to the run containing the "Jeník Hollan, CzechGlobe
#    (..." text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Create the new character styles (wdStyleTypeCharacter = 2)
# ---------------------------------------------------------------------

$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# ---------------------------------------------------------------------
# 2. "Informace v této příručce ..." paragraphs: fix text + apply style
# ---------------------------------------------------------------------

$oldInfo = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od 16. – 25. ledna, 7. – 16. listopadu, 6. – 15. prosince. Při pozorování použijte hvězdy oblohy, které zobrazujíSouhvězdí Perseus.16. – 25. ledna, 7. – 16. listopadu, 6. – 15. prosince"
$newInfo = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od 16. – 25. ledna, 7. – 16. listopadu, 6. – 15. prosince. Při pozorování použijte hvězdy oblohy, které zobrazujíSouhvězdí Perseus."

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Replacement.Style = $ganStyle
$rng.Find.Replacement.Text = $newInfo
$rng.Find.Execute($oldInfo, $true, $false, $false, $false, $false, $true, 1, $false, $newInfo, 2)

# ---------------------------------------------------------------------
# 3. "Jeník Hollan, CzechGlobe (..." run: apply style (text unchanged)
# ---------------------------------------------------------------------

$linkText = "Jeník Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/"

$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found = $rng2.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng2.Style = $ganLinks
}
